# Append additional intraday option-chain rows (Put side) to Sheet1,
# mirroring the existing "Call" rows already present in rows 2-14.
# Columns A-M are stored as text (matching the source data's inline-string
# cells); column N (count) is a genuine number.
#
# Only the columns whose values Excel would otherwise auto-coerce to a
# number/date (strike price, OHLC, volume, open interest, expiry date) get
# an explicit text NumberFormat before the value is written; plain
# alphabetic columns are left with the default format, same as the
# pre-existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("2024-10-31 09:15:00","RELIND","NFO","Options","31-OCT-24","Put","1300",".1",".1",".05",".05","302500","2259000",0),
    @("2024-10-31 09:45:00","RELIND","NFO","Options","31-OCT-24","Put","1300",".1",".1",".05",".1","124000","2176000",1),
    @("2024-10-31 10:15:00","RELIND","NFO","Options","31-OCT-24","Put","1300",".05",".05",".05",".05","21000","2169000",2),
    @("2024-10-31 10:45:00","RELIND","NFO","Options","31-OCT-24","Put","1300",".05",".05",".05",".05","24500","2157500",3),
    @("2024-10-31 11:15:00","RELIND","NFO","Options","31-OCT-24","Put","1300",".05",".05",".05",".05","12500","2133500",4),
    @("2024-10-31 11:45:00","RELIND","NFO","Options","31-OCT-24","Put","1300",".05",".05",".05",".05","2500","2132000",5),
    @("2024-10-31 12:15:00","RELIND","NFO","Options","31-OCT-24","Put","1300",".05",".05",".05",".05","13000","2133000",6),
    @("2024-10-31 12:45:00","RELIND","NFO","Options","31-OCT-24","Put","1300",".05",".05",".05",".05","1000","2133000",7),
    @("2024-10-31 13:15:00","RELIND","NFO","Options","31-OCT-24","Put","1300",".05",".05",".05",".05","6000","2133500",8),
    @("2024-10-31 13:45:00","RELIND","NFO","Options","31-OCT-24","Put","1300",".05",".05",".05",".05","17000","2120500",9),
    @("2024-10-31 14:15:00","RELIND","NFO","Options","31-OCT-24","Put","1300",".05",".05",".05",".05","50000","2099000",10),
    @("2024-10-31 14:45:00","RELIND","NFO","Options","31-OCT-24","Put","1300",".05",".05",".05",".05","1000","2079000",11),
    @("2024-10-31 15:15:00","RELIND","NFO","Options","31-OCT-24","Put","1300",".05",".05",".05",".05","2000","2076000",12)
)

$startRow = 15
# 1-based column indexes that look numeric/date-like and need a text format
# so Excel keeps them as literal strings: E,G,H,I,J,K,L,M
$forceTextCols = @(5,7,8,9,10,11,12,13)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $col = $c + 1
        $cell = $ws.Cells.Item($row, $col)
        if ($col -eq 14) {
            # count: genuine number
            $cell.Value = $values[$c]
        } elseif ($forceTextCols -contains $col) {
            $cell.NumberFormat = "@"
            $cell.Value = [string]$values[$c]
        } else {
            $cell.Value = [string]$values[$c]
        }
    }
}
